# Updates cryptos list values (Price / Volume(1h)) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '40.178.43'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +0.77%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.237.15'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  +0.15%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '293.70'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.71%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '89.01'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +5.67%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -0.05%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '31.22'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +2.44%  '
$ws.Range('E11').Value = '  +1.08%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '47.68'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('E13').Value = '  +1.46%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '6.46'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +1.69%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '2.582.33'
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '14.24'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -0.27%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.212.97'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('E18').Value = '  +2.06%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '40.147.94'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +0.93%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '11.64'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +10.78%  '
$ws.Range('E21').Value = '  +0.84%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '5.87'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('E23').Value = '  +1.16%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '236.68'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +3.02%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.48'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('E27').Value = '  -0.54%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '23.01'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('E30').Value = '  +1.45%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '33.10'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.45%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '152.59'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('E33').Value = '  -0.05%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.00'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('E37').Value = '  +6.99%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '16.35'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  +2.62%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.74'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +2.69%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '2.121.98'
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '3.86'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +2.94%  '
$ws.Range('E44').Value = '  +6.66%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '18.29'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +11.23%  '
$ws.Range('E46').Value = '  +2.63%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '10.13'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +10.55%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.30%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.451.94'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.14%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '71.45'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('E51').Value = '  +5.45%  '
